# Add a new "2022" data column (L) to the freshwater-withdrawal table,
# mirroring the existing 2021 column (K) for formatting, then filling in
# the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clone formatting (number format, font, borders, alignment) from the
#    2021 column (K4:K18) into the new 2022 column (L4:L18) by copying the
#    whole range; values get overwritten individually right after.
$ws.Range("K4:K18").Copy($ws.Range("L4:L18"))

# 2. Header year
$ws.Range("L4").Value = 2022

# 3. Total freshwater withdrawal
$ws.Range("L5").Value = 8800.6

# 4. "by type of source" section header (blank data row)
$ws.Range("L6").Value = $null

# 5. from natural water sources (formula: total - underground)
$ws.Range("L7").Formula = "=L5-L8"

# 6. from underground horizons
$ws.Range("L8").Value = 258.39999999999998

# 7. "by territory" section header (blank data row)
$ws.Range("L9").Value = $null

# 8. Oblasts / cities
$ws.Range("L10").Value = 683.8
$ws.Range("L11").Value = 1101.8
$ws.Range("L12").Value = 714.9
$ws.Range("L13").Value = 757.9
$ws.Range("L14").Value = 1383.3
$ws.Range("L15").Value = 1023.7
$ws.Range("L16").Value = 2929.3
$ws.Range("L17").Value = 148.9
$ws.Range("L18").Value = 57

# 9. Move the active selection to M4 (matches the saved view state in the
#    source workbook after the new column was added).
$ws.Range("M4").Select() | Out-Null

Write-Host "Added 2022 column (L) to sheet1"
